# Auto-generated edit script applying the Ifrit_Profits leve-profit recalculation
# to the 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 3081.1428
$ws.Cells.Item(116, 9).Value = 2336.7273
$ws.Cells.Item(116, 10).Value = 3900
$ws.Cells.Item(116, 11).Value = 2336.7273
$ws.Cells.Item(116, 12).Value = 3900
$ws.Cells.Item(116, 13).Value = 1105.2727
$ws.Cells.Item(116, 14).Value = -10784

$ws.Cells.Item(124, 8).Value = 50980
$ws.Cells.Item(124, 10).Value = 50980
$ws.Cells.Item(124, 12).Value = 50980
$ws.Cells.Item(124, 14).Value = -60800

$ws.Cells.Item(132, 8).Value = 297081.97
$ws.Cells.Item(132, 9).Value = 336562.9
$ws.Cells.Item(132, 11).Value = 1009688.7
$ws.Cells.Item(132, 13).Value = -1007158.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5840.5894
$ws.Cells.Item(32, 9).Value = 6392.5347
$ws.Cells.Item(32, 10).Value = 4014.923
$ws.Cells.Item(32, 11).Value = 6392.5347
$ws.Cells.Item(32, 12).Value = 4014.923
$ws.Cells.Item(32, 13).Value = -6105.5347
$ws.Cells.Item(32, 14).Value = -4588.923

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 38857.734
$ws.Cells.Item(134, 9).Value = 41097.07
$ws.Cells.Item(134, 10).Value = 7507
$ws.Cells.Item(134, 11).Value = 123291.21
$ws.Cells.Item(134, 12).Value = 22521
$ws.Cells.Item(134, 13).Value = -120756.21
$ws.Cells.Item(134, 14).Value = -27591

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1157
$ws.Cells.Item(16, 9).Value = 1151.0834
$ws.Cells.Item(16, 10).Value = 1168.8334
$ws.Cells.Item(16, 11).Value = 1151.0834
$ws.Cells.Item(16, 12).Value = 1168.8334
$ws.Cells.Item(16, 13).Value = -864.0834
$ws.Cells.Item(16, 14).Value = -1742.8334

$ws.Cells.Item(58, 8).Value = 2280.9062
$ws.Cells.Item(58, 9).Value = 1622.75
$ws.Cells.Item(58, 10).Value = 3377.8333
$ws.Cells.Item(58, 11).Value = 1622.75
$ws.Cells.Item(58, 12).Value = 3377.8333
$ws.Cells.Item(58, 13).Value = -1419.75
$ws.Cells.Item(58, 14).Value = -3783.8333

$ws.Cells.Item(86, 8).Value = 5233.8335
$ws.Cells.Item(86, 9).Value = 8100
$ws.Cells.Item(86, 10).Value = 4278.4443
$ws.Cells.Item(86, 11).Value = 8100
$ws.Cells.Item(86, 12).Value = 4278.4443
$ws.Cells.Item(86, 13).Value = -6977
$ws.Cells.Item(86, 14).Value = -6524.4443

$ws.Cells.Item(89, 8).Value = 5233.8335
$ws.Cells.Item(89, 9).Value = 8100
$ws.Cells.Item(89, 10).Value = 4278.4443
$ws.Cells.Item(89, 11).Value = 40500
$ws.Cells.Item(89, 12).Value = 21392.2215
$ws.Cells.Item(89, 13).Value = -34884
$ws.Cells.Item(89, 14).Value = -32624.2215

$ws.Cells.Item(99, 8).Value = 1240.5
$ws.Cells.Item(99, 10).Value = 1225
$ws.Cells.Item(99, 12).Value = 1225
$ws.Cells.Item(99, 14).Value = -4221

$ws.Cells.Item(113, 8).Value = 1157
$ws.Cells.Item(113, 9).Value = 1151.0834
$ws.Cells.Item(113, 10).Value = 1168.8334
$ws.Cells.Item(113, 11).Value = 1151.0834
$ws.Cells.Item(113, 12).Value = 1168.8334
$ws.Cells.Item(113, 13).Value = 1018.9166
$ws.Cells.Item(113, 14).Value = -5508.8334

$ws.Cells.Item(126, 8).Value = 1240.5
$ws.Cells.Item(126, 10).Value = 1225
$ws.Cells.Item(126, 12).Value = 3675
$ws.Cells.Item(126, 14).Value = -8615

$ws.Cells.Item(132, 8).Value = 1751.5
$ws.Cells.Item(132, 9).Value = 1145.9286
$ws.Cells.Item(132, 10).Value = 2811.25
$ws.Cells.Item(132, 11).Value = 3437.7858
$ws.Cells.Item(132, 12).Value = 8433.75
$ws.Cells.Item(132, 13).Value = -907.7857999999997
$ws.Cells.Item(132, 14).Value = -13493.75

$ws.Cells.Item(136, 8).Value = 2280.9062
$ws.Cells.Item(136, 9).Value = 1622.75
$ws.Cells.Item(136, 10).Value = 3377.8333
$ws.Cells.Item(136, 11).Value = 4868.25
$ws.Cells.Item(136, 12).Value = 10133.4999
$ws.Cells.Item(136, 13).Value = -2318.25
$ws.Cells.Item(136, 14).Value = -15233.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 12860.4
$ws.Cells.Item(16, 9).Value = 20100.666
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 60301.99800000001
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = -60128.99800000001
$ws.Cells.Item(16, 14).Value = -6346

$ws.Cells.Item(33, 8).Value = 62500296
$ws.Cells.Item(33, 10).Value = 313.25
$ws.Cells.Item(33, 12).Value = 1879.5
$ws.Cells.Item(33, 14).Value = -2445.5

$ws.Cells.Item(113, 8).Value = 2885.4187
$ws.Cells.Item(113, 9).Value = 530.913
$ws.Cells.Item(113, 10).Value = 5593.1
$ws.Cells.Item(113, 11).Value = 1592.739
$ws.Cells.Item(113, 12).Value = 16779.3
$ws.Cells.Item(113, 13).Value = 577.261
$ws.Cells.Item(113, 14).Value = -21119.3

$ws.Cells.Item(132, 8).Value = 31250902
$ws.Cells.Item(132, 9).Value = 38462224
$ws.Cells.Item(132, 10).Value = 1843.3334
$ws.Cells.Item(132, 11).Value = 346160016
$ws.Cells.Item(132, 12).Value = 16590.0006
$ws.Cells.Item(132, 13).Value = -346157486
$ws.Cells.Item(132, 14).Value = -21650.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 21705.143
$ws.Cells.Item(123, 10).Value = 21705.143
$ws.Cells.Item(123, 12).Value = 21705.143
$ws.Cells.Item(123, 14).Value = -26605.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1443.9
$ws.Cells.Item(40, 9).Value = 1369.0769
$ws.Cells.Item(40, 10).Value = 1582.8572
$ws.Cells.Item(40, 11).Value = 1369.0769
$ws.Cells.Item(40, 12).Value = 1582.8572
$ws.Cells.Item(40, 13).Value = -1233.0769
$ws.Cells.Item(40, 14).Value = -1854.8572

$ws.Cells.Item(61, 8).Value = 1502.3077
$ws.Cells.Item(61, 9).Value = 1093.6364
$ws.Cells.Item(61, 10).Value = 3750
$ws.Cells.Item(61, 11).Value = 1093.6364
$ws.Cells.Item(61, 12).Value = 3750
$ws.Cells.Item(61, 13).Value = -891.6364000000001
$ws.Cells.Item(61, 14).Value = -4154

$ws.Cells.Item(82, 8).Value = 1615.091
$ws.Cells.Item(82, 9).Value = 1375.2
$ws.Cells.Item(82, 11).Value = 1375.2
$ws.Cells.Item(82, 13).Value = -1014.2

$ws.Cells.Item(85, 8).Value = 1615.091
$ws.Cells.Item(85, 9).Value = 1375.2
$ws.Cells.Item(85, 11).Value = 1375.2
$ws.Cells.Item(85, 13).Value = -127.2

$ws.Cells.Item(113, 8).Value = 1502.3077
$ws.Cells.Item(113, 9).Value = 1093.6364
$ws.Cells.Item(113, 10).Value = 3750
$ws.Cells.Item(113, 11).Value = 1093.6364
$ws.Cells.Item(113, 12).Value = 3750
$ws.Cells.Item(113, 13).Value = 1076.3636
$ws.Cells.Item(113, 14).Value = -8090

$ws.Cells.Item(136, 8).Value = 1341.125
$ws.Cells.Item(136, 9).Value = 532.7143
$ws.Cells.Item(136, 11).Value = 1598.1429
$ws.Cells.Item(136, 13).Value = 951.8571000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 598.5
$ws.Cells.Item(113, 9).Value = 536.8889
$ws.Cells.Item(113, 10).Value = 783.3333
$ws.Cells.Item(113, 11).Value = 1610.6667
$ws.Cells.Item(113, 12).Value = 2349.9999
$ws.Cells.Item(113, 13).Value = 559.3332999999998
$ws.Cells.Item(113, 14).Value = -6689.9999

$ws.Cells.Item(132, 8).Value = 3129.4565
$ws.Cells.Item(132, 9).Value = 3543
$ws.Cells.Item(132, 11).Value = 10629
$ws.Cells.Item(132, 13).Value = -8099

$ws.Cells.Item(136, 8).Value = 9876.682000000001
$ws.Cells.Item(136, 9).Value = 10204.143
$ws.Cells.Item(136, 11).Value = 30612.429
$ws.Cells.Item(136, 13).Value = -28062.429
